# Atualizado por script em 01-11-2023 20:45
#
# 1) Swap match data (columns F:V) between rows 12 and 13.
# 2) Rotate match data (columns F:V) across rows 37 -> 38 -> 40 -> 37.
# 3) Append a new match result as row 74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap rows 12 and 13 (F:V) ---
$row12 = $ws.Range("F12:V12").Value2
$row13 = $ws.Range("F13:V13").Value2
$ws.Range("F12:V12").Value2 = $row13
$ws.Range("F13:V13").Value2 = $row12

# --- 2) Rotate rows 37, 38, 40 (F:V): new37=old38, new38=old40, new40=old37 ---
$row37 = $ws.Range("F37:V37").Value2
$row38 = $ws.Range("F38:V38").Value2
$row40 = $ws.Range("F40:V40").Value2
$ws.Range("F37:V37").Value2 = $row38
$ws.Range("F38:V38").Value2 = $row40
$ws.Range("F40:V40").Value2 = $row37

# --- 3) Append new row 74, mirroring the formatting of the last existing row (73) ---
$ws.Range("A73:V73").Copy()
$ws.Range("A74:V74").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(74, 1).Value2 = 73
$ws.Cells.Item(74, 2).Value2 = "south-africa"
$ws.Cells.Item(74, 3).Value2 = "premier-league"
$ws.Cells.Item(74, 4).Value2 = "2023-2024"
$ws.Cells.Item(74, 5).Value2 = 45231.77083333334
$ws.Cells.Item(74, 6).Value2 = "Cape Town Spurs"
$ws.Cells.Item(74, 7).Value2 = 2
$ws.Cells.Item(74, 8).Value2 = "Orlando Pirates"
$ws.Cells.Item(74, 9).Value2 = 1
$ws.Cells.Item(74, 10).Value2 = 5.93
$ws.Cells.Item(74, 11).Value2 = "01/11/2023 13:42"
$ws.Cells.Item(74, 12).Value2 = 6.03
$ws.Cells.Item(74, 13).Value2 = "01/11/2023 18:29"
$ws.Cells.Item(74, 14).Value2 = 3.97
$ws.Cells.Item(74, 15).Value2 = "01/11/2023 13:42"
$ws.Cells.Item(74, 16).Value2 = 3.81
$ws.Cells.Item(74, 17).Value2 = "01/11/2023 18:29"
$ws.Cells.Item(74, 18).Value2 = 1.56
$ws.Cells.Item(74, 19).Value2 = "01/11/2023 13:42"
$ws.Cells.Item(74, 20).Value2 = 1.61
$ws.Cells.Item(74, 21).Value2 = "01/11/2023 18:29"
$ws.Cells.Item(74, 22).Value2 = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-orlando-pirates/ngpe5dHP/"
